$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("INS")

# Update the TIMESTEP command from 15 to 12
$ws.Range("D9").Value = '$SET TIMESTEP 12'

# Clear any explicit (duplicate) cell style on B10 so it reverts to the default style
$ws.Range("B10").Style = "Normal"

# Update the active selection to match the edited cell
$ws.Range("D9").Select()
